$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 148.34
$ws.Range("I15").Value = 148.34
$ws.Range("K15").Value = 445.02
$ws.Range("M15").Value = -276.02
$ws.Range("H39").Value = 393.13043
$ws.Range("I39").Value = 51.2
$ws.Range("J39").Value = 1034.25
$ws.Range("K39").Value = 153.6
$ws.Range("L39").Value = 3102.75
$ws.Range("M39").Value = 142.4
$ws.Range("N39").Value = -3694.75
$ws.Range("H62").Value = 2222.2666
$ws.Range("I62").Value = 1343.4
$ws.Range("J62").Value = 3980
$ws.Range("K62").Value = 1343.4
$ws.Range("L62").Value = 3980
$ws.Range("M62").Value = -719.4000000000001
$ws.Range("N62").Value = -5228
$ws.Range("H65").Value = 2222.2666
$ws.Range("I65").Value = 1343.4
$ws.Range("J65").Value = 3980
$ws.Range("K65").Value = 6717
$ws.Range("L65").Value = 19900
$ws.Range("M65").Value = -3597
$ws.Range("N65").Value = -26140
$ws.Range("H82").Value = 2115
$ws.Range("I82").Value = 655
$ws.Range("K82").Value = 1965
$ws.Range("M82").Value = -1559
$ws.Range("H85").Value = 2115
$ws.Range("I85").Value = 655
$ws.Range("K85").Value = 1965
$ws.Range("M85").Value = -561
$ws.Range("H88").Value = 2288.9697
$ws.Range("J88").Value = 1251.75
$ws.Range("L88").Value = 1251.75
$ws.Range("N88").Value = -2063.75
$ws.Range("H91").Value = 2288.9697
$ws.Range("J91").Value = 1251.75
$ws.Range("L91").Value = 1251.75
$ws.Range("N91").Value = -4059.75
$ws.Range("H98").Value = 752.43475
$ws.Range("I98").Value = 773.3333
$ws.Range("J98").Value = 533
$ws.Range("K98").Value = 773.3333
$ws.Range("L98").Value = 533
$ws.Range("M98").Value = 724.6667
$ws.Range("N98").Value = -3529
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("M119").ClearContents()
$ws.Range("H122").Value = 752.43475
$ws.Range("I122").Value = 773.3333
$ws.Range("J122").Value = 533
$ws.Range("K122").Value = 2319.9999
$ws.Range("L122").Value = 1599
$ws.Range("M122").Value = 130.0001000000002
$ws.Range("N122").Value = -6499
$ws.Range("H137").Value = 56861.5
$ws.Range("I137").Value = 999.5833
$ws.Range("K137").Value = 2998.7499
$ws.Range("M137").Value = -448.7498999999998
$ws.Range("H138").Value = 5185.18
$ws.Range("I138").Value = 1453.2812
$ws.Range("J138").Value = 6941.3677
$ws.Range("K138").Value = 4359.8436
$ws.Range("L138").Value = 20824.1031
$ws.Range("M138").Value = 780.1563999999998
$ws.Range("N138").Value = -31104.1031
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 179.53847
$ws.Range("I80").Value = 15
$ws.Range("J80").Value = 193.25
$ws.Range("K80").Value = 15
$ws.Range("L80").Value = 193.25
$ws.Range("M80").Value = 983
$ws.Range("N80").Value = -2189.25
$ws.Range("H83").Value = 179.53847
$ws.Range("I83").Value = 15
$ws.Range("J83").Value = 193.25
$ws.Range("K83").Value = 75
$ws.Range("L83").Value = 966.25
$ws.Range("M83").Value = 4917
$ws.Range("N83").Value = -10950.25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 866
$ws.Range("I105").Value = 806.3333
$ws.Range("J105").Value = 955.5
$ws.Range("K105").Value = 806.3333
$ws.Range("L105").Value = 955.5
$ws.Range("M105").Value = 940.6667
$ws.Range("N105").Value = -4449.5
$ws.Range("H132").Value = 1701.1143
$ws.Range("I132").Value = 996.0454999999999
$ws.Range("J132").Value = 2894.3076
$ws.Range("K132").Value = 2988.1365
$ws.Range("L132").Value = 8682.9228
$ws.Range("M132").Value = -458.1364999999996
$ws.Range("N132").Value = -13742.9228
$ws.Range("H134").Value = 2183.0833
$ws.Range("I134").Value = 2270.8445
$ws.Range("J134").Value = 866.6667
$ws.Range("K134").Value = 6812.5335
$ws.Range("L134").Value = 2600.0001
$ws.Range("M134").Value = -4277.5335
$ws.Range("N134").Value = -7670.0001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 735.39655
$ws.Range("I5").Value = 538.95557
$ws.Range("J5").Value = 1415.3846
$ws.Range("K5").Value = 1616.86671
$ws.Range("L5").Value = 4246.1538
$ws.Range("M5").Value = -1504.86671
$ws.Range("N5").Value = -4470.1538
$ws.Range("H12").Value = 62.703705
$ws.Range("I12").Value = 8.199999999999999
$ws.Range("J12").Value = 94.76470999999999
$ws.Range("K12").Value = 24.6
$ws.Range("L12").Value = 284.29413
$ws.Range("M12").Value = 148.4
$ws.Range("N12").Value = -630.29413
$ws.Range("H17").Value = 575
$ws.Range("I17").Value = 400
$ws.Range("J17").Value = 633.3333
$ws.Range("K17").Value = 1200
$ws.Range("L17").Value = 1899.9999
$ws.Range("M17").Value = -1031
$ws.Range("N17").Value = -2237.9999
$ws.Range("H33").Value = 5813.8335
$ws.Range("I33").Value = 218.875
$ws.Range("J33").Value = 10289.8
$ws.Range("K33").Value = 1313.25
$ws.Range("L33").Value = 61738.8
$ws.Range("M33").Value = -1030.25
$ws.Range("N33").Value = -62304.8
$ws.Range("H34").Value = 356.9394
$ws.Range("I34").Value = 160.75
$ws.Range("J34").Value = 469.0476
$ws.Range("K34").Value = 482.25
$ws.Range("L34").Value = 1407.1428
$ws.Range("M34").Value = -398.25
$ws.Range("N34").Value = -1575.1428
$ws.Range("H51").Value = 3013.1904
$ws.Range("I51").Value = 484.85715
$ws.Range("J51").Value = 4277.357
$ws.Range("K51").Value = 1454.57145
$ws.Range("L51").Value = 12832.071
$ws.Range("M51").Value = -994.5714499999999
$ws.Range("N51").Value = -13752.071
$ws.Range("H55").Value = 2320
$ws.Range("J55").Value = 2467.2727
$ws.Range("L55").Value = 7401.8181
$ws.Range("N55").Value = -7755.8181
$ws.Range("H64").Value = 4800.28
$ws.Range("I64").Value = 762
$ws.Range("J64").Value = 5151.4346
$ws.Range("K64").Value = 2286
$ws.Range("L64").Value = 15454.3038
$ws.Range("M64").Value = -2016
$ws.Range("N64").Value = -15994.3038
$ws.Range("H67").Value = 4800.28
$ws.Range("I67").Value = 762
$ws.Range("J67").Value = 5151.4346
$ws.Range("K67").Value = 2286
$ws.Range("L67").Value = 15454.3038
$ws.Range("M67").Value = -1350
$ws.Range("N67").Value = -17326.3038
$ws.Range("H98").Value = 658.5
$ws.Range("I98").Value = 850
$ws.Range("J98").Value = 467
$ws.Range("K98").Value = 2550
$ws.Range("L98").Value = 1401
$ws.Range("M98").Value = -1052
$ws.Range("N98").Value = -4397
$ws.Range("H113").Value = 2006.5
$ws.Range("I113").Value = 498
$ws.Range("J113").Value = 3515
$ws.Range("K113").Value = 1494
$ws.Range("L113").Value = 10545
$ws.Range("M113").Value = 676
$ws.Range("N113").Value = -14885
$ws.Range("H135").Value = 735.39655
$ws.Range("I135").Value = 538.95557
$ws.Range("J135").Value = 1415.3846
$ws.Range("K135").Value = 4850.60013
$ws.Range("L135").Value = 12738.4614
$ws.Range("M135").Value = -2315.60013
$ws.Range("N135").Value = -17808.4614
$ws.Range("H140").Value = 1593.2142
$ws.Range("I140").Value = 1091.7391
$ws.Range("J140").Value = 3900
$ws.Range("K140").Value = 3275.2173
$ws.Range("L140").Value = 11700
$ws.Range("M140").Value = 1904.7827
$ws.Range("N140").Value = -22060
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 3072.3333
$ws.Range("I43").Value = 1686.8
$ws.Range("J43").Value = 10000
$ws.Range("K43").Value = 1686.8
$ws.Range("L43").Value = 10000
$ws.Range("M43").Value = -1535.8
$ws.Range("N43").Value = -10302
$ws.Range("H122").Value = 2209.5833
$ws.Range("I122").Value = 2274.0908
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 6822.2724
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -4372.2724
$ws.Range("N122").Value = -9400
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1667473.6
$ws.Range("I22").Value = 2778342.8
$ws.Range("J22").Value = 1170
$ws.Range("K22").Value = 2778342.8
$ws.Range("L22").Value = 1170
$ws.Range("M22").Value = -2778047.8
$ws.Range("N22").Value = -1760
$ws.Range("H27").Value = 1667473.6
$ws.Range("I27").Value = 2778342.8
$ws.Range("J27").Value = 1170
$ws.Range("K27").Value = 2778342.8
$ws.Range("L27").Value = 1170
$ws.Range("M27").Value = -2778235.8
$ws.Range("N27").Value = -1384
$ws.Range("H122").Value = 7043.278
$ws.Range("I122").Value = 9998.9
$ws.Range("J122").Value = 3348.75
$ws.Range("K122").Value = 29996.7
$ws.Range("L122").Value = 10046.25
$ws.Range("M122").Value = -27546.7
$ws.Range("N122").Value = -14946.25
$ws.Range("H132").Value = 11634918
$ws.Range("I132").Value = 23821074
$ws.Range("J132").Value = 2678.5908
$ws.Range("K132").Value = 71463222
$ws.Range("L132").Value = 8035.7724
$ws.Range("M132").Value = -71460692
$ws.Range("N132").Value = -13095.7724
$ws.Range("H136").Value = 9682.611000000001
$ws.Range("I136").Value = 18786.857
$ws.Range("J136").Value = 3889
$ws.Range("K136").Value = 56360.571
$ws.Range("L136").Value = 11667
$ws.Range("M136").Value = -53810.571
$ws.Range("N136").Value = -16767
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 52633724
$ws.Range("I81").Value = 100001920
$ws.Range("J81").Value = 2387.7778
$ws.Range("K81").Value = 200003840
$ws.Range("L81").Value = 4775.5556
$ws.Range("M81").Value = -200002779
$ws.Range("N81").Value = -6897.5556
$ws.Range("H84").Value = 52633724
$ws.Range("I84").Value = 100001920
$ws.Range("J84").Value = 2387.7778
$ws.Range("K84").Value = 1000019200
$ws.Range("L84").Value = 23877.778
$ws.Range("M84").Value = -1000013896
$ws.Range("N84").Value = -34485.778
